$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to be treated as text so that values
# like "1.002" or "0.05440" are not auto-converted to numbers,
# which would strip significant trailing/leading characters.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "28.680.56"
$ws.Range("E2").Value = "  -1.87%  "

# Row 3
$ws.Range("D3").Value = "1.802.48"
$ws.Range("E3").Value = "  -1.47%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").Value = "231.26"
$ws.Range("E5").Value = "  -2.38%  "

# Row 6
$ws.Range("D6").Value = "0.5948"
$ws.Range("E6").Value = "  -2.23%  "

# Row 7
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.21%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.2769"
$ws.Range("E8").Value = "  -1.81%  "

# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.06807"
$ws.Range("E9").Value = "  -4.36%  "

# Row 10
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "23.29"
$ws.Range("E10").Value = "  -2.74%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07509"
$ws.Range("E11").Value = "  -2.22%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.806.97"
$ws.Range("E12").Value = "  -1.09%  "

# Row 13
$ws.Range("D13").Value = "4.690"
$ws.Range("E13").Value = "  -2.64%  "

# Row 14
$ws.Range("D14").Value = "0.6283"
$ws.Range("E14").Value = "  -1.26%  "

# Row 15
$ws.Range("D15").Value = "2.048.25"
$ws.Range("E15").Value = "  -1.41%  "

# Row 16
$ws.Range("D16").Value = "0.000009202"
$ws.Range("E16").Value = "  -8.05%  "

# Row 17
$ws.Range("D17").Value = "75.38"
$ws.Range("E17").Value = "  -5.10%  "

# Row 18
$ws.Range("D18").Value = "28.566.06"
$ws.Range("E18").Value = "  -2.19%  "

# Row 19
$ws.Range("D19").Value = "5.465"
$ws.Range("E19").Value = "  -7.45%  "

# Row 20
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.16%  "

# Row 21
$ws.Range("D21").Value = "209.05"
$ws.Range("E21").Value = "  -8.83%  "

# Row 22
$ws.Range("D22").Value = "11.37"
$ws.Range("E22").Value = "  -3.90%  "

# Row 23
$ws.Range("D23").Value = "6.804"
$ws.Range("E23").Value = "  -3.18%  "

# Row 24
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.29%  "

# Row 25
$ws.Range("D25").Value = "154.74"
$ws.Range("E25").Value = "  +0.10%  "

# Row 26
$ws.Range("D26").Value = "7.849"
$ws.Range("E26").Value = "  -2.99%  "

# Row 27
$ws.Range("D27").Value = "0.1273"
$ws.Range("E27").Value = "  -1.34%  "

# Row 28
$ws.Range("D28").Value = "16.38"
$ws.Range("E28").Value = "  -1.68%  "

# Row 29
$ws.Range("D29").Value = "1.446"
$ws.Range("E29").Value = "  -4.37%  "

# Row 30
$ws.Range("D30").Value = "0.06319"
$ws.Range("E30").Value = "  -3.60%  "

# Row 31
$ws.Range("D31").Value = "1.418"
$ws.Range("E31").Value = "  -2.52%  "

# Row 32
$ws.Range("D32").Value = "3.740"
$ws.Range("E32").Value = "  -2.79%  "

# Row 33
$ws.Range("D33").Value = "3.721"
$ws.Range("E33").Value = "  -3.20%  "

# Row 34
$ws.Range("D34").Value = "1.719"
$ws.Range("E34").Value = "  -1.07%  "

# Row 35
$ws.Range("D35").Value = "1.049"
$ws.Range("E35").Value = "  -7.48%  "

# Row 36
$ws.Range("D36").Value = "0.6356"
$ws.Range("E36").Value = "  -2.93%  "

# Row 37
$ws.Range("D37").Value = "2.506"
$ws.Range("E37").Value = "  -1.48%  "

# Row 38
$ws.Range("D38").Value = "2.720"
$ws.Range("E38").Value = "  -1.72%  "

# Row 39
$ws.Range("D39").Value = "0.01701"
$ws.Range("E39").Value = "  -3.41%  "

# Row 40
$ws.Range("D40").Value = "6.404"
$ws.Range("E40").Value = "  -3.45%  "

# Row 41
$ws.Range("D41").Value = "1.137.22"
$ws.Range("E41").Value = "  -6.93%  "

# Row 42
$ws.Range("D42").Value = "0.8590"
$ws.Range("E42").Value = "  -7.48%  "

# Row 43
$ws.Range("D43").Value = "1.003"
$ws.Range("E43").Value = "  +0.25%  "

# Row 44
$ws.Range("D44").Value = "100.96"
$ws.Range("E44").Value = "  -0.13%  "

# Row 45
$ws.Range("D45").Value = "1.960.33"
$ws.Range("E45").Value = "  -1.05%  "

# Row 46
$ws.Range("D46").Value = "60.54"
$ws.Range("E46").Value = "  -4.82%  "

# Row 47
$ws.Range("D47").Value = "0.00000000111"
$ws.Range("E47").Value = "  -4.99%  "

# Row 48
$ws.Range("D48").Value = "1.580"
$ws.Range("E48").Value = "  -2.01%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.319"
$ws.Range("E49").Value = "  -2.98%  "

# Row 50
$ws.Range("D50").Value = "0.4498"
$ws.Range("E50").Value = "  -1.29%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05440"
$ws.Range("E51").Value = "  -1.83%  "
